$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-07"

# Update the shared-string label used in A9 (August row header)
$ws.Range("A9").Value = "August (through 08-07)"

# Update the August row (row 9) values for 2015-2022
$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 9
$ws.Range("D9").Value = 20
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 12
$ws.Range("G9").Value = 38
$ws.Range("H9").Value = 42
$ws.Range("I9").Value = 38

# Update the Total row (row 10) values for 2015-2022
$ws.Range("B10").Value = 171
$ws.Range("C10").Value = 311
$ws.Range("D10").Value = 485
$ws.Range("E10").Value = 440
$ws.Range("F10").Value = 316
$ws.Range("G10").Value = 659
$ws.Range("H10").Value = 952
$ws.Range("I10").Value = 1008
